$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card7")

# Header M1: drop the trailing space ("Event " -> "Event")
$ws.Cells.Item(1, 13).Value = "Event"

# New header N1: "Correction " (keep trailing space), same header style as the rest of row 1
$ws.Cells.Item(1, 14).Value = "Correction "
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2..12: M gets the "nan" placeholder text (matching the other data columns),
# N is added as a new (empty) column so the used range grows to A1:N12.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 13).Value = "nan"

    $ws.Cells.Item($r, 14).Value = ""
    $ws.Cells.Item($r, 14).Style = $ws.Cells.Item($r, 12).Style
}
